# question3 with excel activities
$wb = $excel.ActiveWorkbook

# Add a new worksheet and move it before the existing "Sum" sheet so it
# becomes the first tab in the workbook.
$sumSheet = $wb.Worksheets.Item("Sum")
$newSheet = $wb.Worksheets.Add()
$newSheet.Move($sumSheet)
$newSheet.Name = "Sheet1"

# Header row
$newSheet.Range("A1").Value = "CashIn"
$newSheet.Range("B1").Value = "OnCheck"
$newSheet.Range("C1").Value = "NotOnCheck"

# Data row
$newSheet.Range("A2").Value = 35630
$newSheet.Range("B2").Value = 179809
$newSheet.Range("C2").Value = 19565

$newSheet.Activate()
